$d = $word.ActiveDocument
$ellipsis = [char]0x2026

# --- 1) "…" -> "1" ----------------------------------------------------------
# The caption paragraph reads "ตาราง … Activity Diagram" (Table … Activity
# Diagram). Replace the ellipsis character itself with "1", leaving the
# space that already follows it untouched.
$rNum = $d.Content
$foundNum = $rNum.Find.Execute($ellipsis, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundNum) {
    throw "could not find the ellipsis character"
}
$rNum.Text = "1"

# --- 2) lone space after "ตาราง" -> "ที่ " -----------------------------------
# Turn "ตาราง " into "ตารางที่ " (i.e. "Table" -> "Table No.") by rewriting
# just the single space run that follows "ตาราง".
$rLabel = $d.Content
$foundLabel = $rLabel.Find.Execute("ตาราง ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundLabel) {
    throw "could not find 'ตาราง '"
}
$rSpace = $d.Range($rLabel.End - 1, $rLabel.End)
$rSpace.Text = "ที่ "
